# Auto-generated edit script: update TPM-derived metrics for Fgf2-Fgfr4 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.7321483333333333"
$ws.Range("H2").Value = [double]"2.196445"
$ws.Range("I2").Value = [double]"0.05113520435363902"
$ws.Range("J2").Value = [double]"0.05113520435363902"
$ws.Range("M2").Value = [double]"0.789222"
$ws.Range("N2").Value = [double]"2.367666"
$ws.Range("O2").Value = [double]"0.01341929863527565"
$ws.Range("P2").Value = [double]"0.01341929863527565"
$ws.Range("Q2").Value = [double]"0.57782757193"
$ws.Range("R2").Value = [double]"5.200448147369999"
$ws.Range("S2").Value = [double]"0.0006861985779973294"
$ws.Range("T2").Value = [double]"0.0006861985779973294"
$ws.Range("G3").Value = [double]"0.7321483333333333"
$ws.Range("H3").Value = [double]"2.196445"
$ws.Range("I3").Value = [double]"0.05113520435363902"
$ws.Range("J3").Value = [double]"0.05113520435363902"
$ws.Range("O3").Value = [double]"0.005047365584441773"
$ws.Range("P3").Value = [double]"0.005047365584441773"
$ws.Range("Q3").Value = [double]"0.2173367684533333"
$ws.Range("R3").Value = [double]"1.95603091608"
$ws.Range("S3").Value = [double]"0.0002580980706079547"
$ws.Range("T3").Value = [double]"0.0002580980706079547"
$ws.Range("G4").Value = [double]"0.7321483333333333"
$ws.Range("H4").Value = [double]"2.196445"
$ws.Range("I4").Value = [double]"0.05113520435363902"
$ws.Range("J4").Value = [double]"0.05113520435363902"
$ws.Range("M4").Value = [double]"57.61405833333333"
$ws.Range("N4").Value = [double]"172.842175"
$ws.Range("O4").Value = [double]"0.9796232927683105"
$ws.Range("P4").Value = [double]"0.9796232927683105"
$ws.Range("Q4").Value = [double]"42.18203678531944"
$ws.Range("R4").Value = [double]"379.638331067875"
$ws.Range("S4").Value = [double]"0.0500932372652923"
$ws.Range("T4").Value = [double]"0.0500932372652923"
$ws.Range("G5").Value = [double]"0.7321483333333333"
$ws.Range("H5").Value = [double]"2.196445"
$ws.Range("I5").Value = [double]"0.05113520435363902"
$ws.Range("J5").Value = [double]"0.05113520435363902"
$ws.Range("M5").Value = [double]"0.1123343333333333"
$ws.Range("N5").Value = [double]"0.337003"
$ws.Range("O5").Value = [double]"0.001910043011972043"
$ws.Range("P5").Value = [double]"0.001910043011972043"
$ws.Range("Q5").Value = [double]"0.0822453949261111"
$ws.Range("R5").Value = [double]"0.7402085543349999"
$ws.Range("S5").Value = [double]"9.767043974143058E-05"
$ws.Range("T5").Value = [double]"9.76704397414306E-05"
$ws.Range("I6").Value = [double]"0.7165747117895102"
$ws.Range("J6").Value = [double]"0.7165747117895102"
$ws.Range("M6").Value = [double]"0.789222"
$ws.Range("N6").Value = [double]"2.367666"
$ws.Range("O6").Value = [double]"0.01341929863527565"
$ws.Range("P6").Value = [double]"0.01341929863527565"
$ws.Range("Q6").Value = [double]"8.097290918332"
$ws.Range("R6").Value = [double]"72.87561826498799"
$ws.Range("S6").Value = [double]"0.009615930051990015"
$ws.Range("T6").Value = [double]"0.009615930051990015"
$ws.Range("I7").Value = [double]"0.7165747117895102"
$ws.Range("J7").Value = [double]"0.7165747117895102"
$ws.Range("O7").Value = [double]"0.005047365584441773"
$ws.Range("P7").Value = [double]"0.005047365584441773"
$ws.Range("S7").Value = [double]"0.003616814538967656"
$ws.Range("T7").Value = [double]"0.003616814538967656"
$ws.Range("I8").Value = [double]"0.7165747117895102"
$ws.Range("J8").Value = [double]"0.7165747117895102"
$ws.Range("M8").Value = [double]"57.61405833333333"
$ws.Range("N8").Value = [double]"172.842175"
$ws.Range("O8").Value = [double]"0.9796232927683105"
$ws.Range("P8").Value = [double]"0.9796232927683105"
$ws.Range("Q8").Value = [double]"591.1109818412945"
$ws.Range("R8").Value = [double]"5319.99883657165"
$ws.Range("S8").Value = [double]"0.7019732786777431"
$ws.Range("T8").Value = [double]"0.7019732786777431"
$ws.Range("I9").Value = [double]"0.7165747117895102"
$ws.Range("J9").Value = [double]"0.7165747117895102"
$ws.Range("M9").Value = [double]"0.1123343333333333"
$ws.Range("N9").Value = [double]"0.337003"
$ws.Range("O9").Value = [double]"0.001910043011972043"
$ws.Range("P9").Value = [double]"0.001910043011972043"
$ws.Range("Q9").Value = [double]"1.152532211617111"
$ws.Range("R9").Value = [double]"10.372789904554"
$ws.Range("S9").Value = [double]"0.001368688520809434"
$ws.Range("T9").Value = [double]"0.001368688520809435"
$ws.Range("G10").Value = [double]"2.568000333333333"
$ws.Range("H10").Value = [double]"7.704001"
$ws.Range("I10").Value = [double]"0.1793560346266988"
$ws.Range("J10").Value = [double]"0.1793560346266988"
$ws.Range("M10").Value = [double]"0.789222"
$ws.Range("N10").Value = [double]"2.367666"
$ws.Range("O10").Value = [double]"0.01341929863527565"
$ws.Range("P10").Value = [double]"0.01341929863527565"
$ws.Range("Q10").Value = [double]"2.026722359074"
$ws.Range("R10").Value = [double]"18.240501231666"
$ws.Range("S10").Value = [double]"0.002406832190694511"
$ws.Range("T10").Value = [double]"0.002406832190694511"
$ws.Range("G11").Value = [double]"2.568000333333333"
$ws.Range("H11").Value = [double]"7.704001"
$ws.Range("I11").Value = [double]"0.1793560346266988"
$ws.Range("J11").Value = [double]"0.1793560346266988"
$ws.Range("O11").Value = [double]"0.005047365584441773"
$ws.Range("P11").Value = [double]"0.005047365584441773"
$ws.Range("Q11").Value = [double]"0.7623057629493334"
$ws.Range("R11").Value = [double]"6.860751866544"
$ws.Range("S11").Value = [double]"0.0009052754765367465"
$ws.Range("T11").Value = [double]"0.0009052754765367465"
$ws.Range("G12").Value = [double]"2.568000333333333"
$ws.Range("H12").Value = [double]"7.704001"
$ws.Range("I12").Value = [double]"0.1793560346266988"
$ws.Range("J12").Value = [double]"0.1793560346266988"
$ws.Range("M12").Value = [double]"57.61405833333333"
$ws.Range("N12").Value = [double]"172.842175"
$ws.Range("O12").Value = [double]"0.9796232927683105"
$ws.Range("P12").Value = [double]"0.9796232927683105"
$ws.Range("Q12").Value = [double]"147.9529210046861"
$ws.Range("R12").Value = [double]"1331.576289042175"
$ws.Range("S12").Value = [double]"0.1757013492188738"
$ws.Range("T12").Value = [double]"0.1757013492188738"
$ws.Range("G13").Value = [double]"2.568000333333333"
$ws.Range("H13").Value = [double]"7.704001"
$ws.Range("I13").Value = [double]"0.1793560346266988"
$ws.Range("J13").Value = [double]"0.1793560346266988"
$ws.Range("M13").Value = [double]"0.1123343333333333"
$ws.Range("N13").Value = [double]"0.337003"
$ws.Range("O13").Value = [double]"0.001910043011972043"
$ws.Range("P13").Value = [double]"0.001910043011972043"
$ws.Range("Q13").Value = [double]"0.2884746054447778"
$ws.Range("R13").Value = [double]"2.596271449003"
$ws.Range("S13").Value = [double]"0.0003425777405937418"
$ws.Range("T13").Value = [double]"0.0003425777405937418"
$ws.Range("G14").Value = [double]"0.7579039999999999"
$ws.Range("H14").Value = [double]"2.273712"
$ws.Range("I14").Value = [double]"0.05293404923015203"
$ws.Range("J14").Value = [double]"0.05293404923015203"
$ws.Range("M14").Value = [double]"0.789222"
$ws.Range("N14").Value = [double]"2.367666"
$ws.Range("O14").Value = [double]"0.01341929863527565"
$ws.Range("P14").Value = [double]"0.01341929863527565"
$ws.Range("Q14").Value = [double]"0.598154510688"
$ws.Range("R14").Value = [double]"5.383390596191999"
$ws.Range("S14").Value = [double]"0.000710337814593793"
$ws.Range("T14").Value = [double]"0.000710337814593793"
$ws.Range("G15").Value = [double]"0.7579039999999999"
$ws.Range("H15").Value = [double]"2.273712"
$ws.Range("I15").Value = [double]"0.05293404923015203"
$ws.Range("J15").Value = [double]"0.05293404923015203"
$ws.Range("O15").Value = [double]"0.005047365584441773"
$ws.Range("P15").Value = [double]"0.005047365584441773"
$ws.Range("Q15").Value = [double]"0.224982286592"
$ws.Range("R15").Value = [double]"2.024840579328"
$ws.Range("S15").Value = [double]"0.0002671774983294159"
$ws.Range("T15").Value = [double]"0.0002671774983294159"
$ws.Range("G16").Value = [double]"0.7579039999999999"
$ws.Range("H16").Value = [double]"2.273712"
$ws.Range("I16").Value = [double]"0.05293404923015203"
$ws.Range("J16").Value = [double]"0.05293404923015203"
$ws.Range("M16").Value = [double]"57.61405833333333"
$ws.Range("N16").Value = [double]"172.842175"
$ws.Range("O16").Value = [double]"0.9796232927683105"
$ws.Range("P16").Value = [double]"0.9796232927683105"
$ws.Range("Q16").Value = [double]"43.66592526706666"
$ws.Range("R16").Value = [double]"392.9933274035999"
$ws.Range("S16").Value = [double]"0.05185542760640138"
$ws.Range("T16").Value = [double]"0.05185542760640138"
$ws.Range("G17").Value = [double]"0.7579039999999999"
$ws.Range("H17").Value = [double]"2.273712"
$ws.Range("I17").Value = [double]"0.05293404923015203"
$ws.Range("J17").Value = [double]"0.05293404923015203"
$ws.Range("M17").Value = [double]"0.1123343333333333"
$ws.Range("N17").Value = [double]"0.337003"
$ws.Range("O17").Value = [double]"0.001910043011972043"
$ws.Range("P17").Value = [double]"0.001910043011972043"
$ws.Range("Q17").Value = [double]"0.08513864057066665"
$ws.Range("R17").Value = [double]"0.7662477651359999"
$ws.Range("S17").Value = [double]"0.000101106310827436"
$ws.Range("T17").Value = [double]"0.000101106310827436"
